$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    [PSCustomObject]@{ Row=2; D=44230; J=500; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=3; D=44237; J=600; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=4; D=44600; J=400; K=30000; L=31000; M=30500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1220; Q=25 }
    [PSCustomObject]@{ Row=5; D=44225; J=600; K=31000; L=32000; M=31500; N="`$/malla 25 kilos"; O="Provincia de Limarí"; P=1260; Q=25 }
    [PSCustomObject]@{ Row=6; D=44610; J=400; K=30000; L=31000; M=30500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1220; Q=25 }
    [PSCustomObject]@{ Row=7; D=44277; J=560; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=8; D=44663; J=560; K=24000; L=25000; M=24500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=980; Q=25 }
    [PSCustomObject]@{ Row=9; D=44615; J=520; K=31000; L=32000; M=31500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1260; Q=25 }
    [PSCustomObject]@{ Row=10; D=44323; J=600; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=11; D=44655; J=440; K=27000; L=28000; M=27500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1100; Q=25 }
    [PSCustomObject]@{ Row=12; D=44608; J=500; K=30000; L=31000; M=30500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1220; Q=25 }
    [PSCustomObject]@{ Row=13; D=44575; J=400; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=14; D=44547; J=200; K=22000; L=22500; M=22250; N="`$/caja 15 kilos"; O="Provincia de Limarí"; P=1483; Q=15 }
    [PSCustomObject]@{ Row=15; D=44607; J=600; K=30000; L=31000; M=30500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1220; Q=25 }
    [PSCustomObject]@{ Row=16; D=44648; J=480; K=27000; L=28000; M=27500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1100; Q=25 }
    [PSCustomObject]@{ Row=17; D=44662; J=460; K=24000; L=25000; M=24500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=980; Q=25 }
    [PSCustomObject]@{ Row=18; D=44614; J=540; K=31000; L=32000; M=31500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1260; Q=25 }
    [PSCustomObject]@{ Row=19; D=44202; J=600; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia de Limarí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=20; D=44685; J=440; K=26000; L=27000; M=26500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1060; Q=25 }
    [PSCustomObject]@{ Row=21; D=44272; J=600; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=22; D=44692; J=400; K=25000; L=26000; M=25500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1020; Q=25 }
    [PSCustomObject]@{ Row=23; D=44557; J=400; K=30000; L=31000; M=30500; N="`$/malla 25 kilos"; O="Provincia de Limarí"; P=1220; Q=25 }
    [PSCustomObject]@{ Row=24; D=44638; J=400; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=25; D=44274; J=600; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=26; D=44552; J=400; K=35000; L=36000; M=35500; N="`$/malla 25 kilos"; O="Provincia de Limarí"; P=1420; Q=25 }
    [PSCustomObject]@{ Row=27; D=44642; J=400; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=28; D=44291; J=500; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=29; D=44252; J=520; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=30; D=44566; J=400; K=16000; L=17000; M=16500; N="`$/caja 15 kilos"; O="Provincia de Limarí"; P=1100; Q=15 }
    [PSCustomObject]@{ Row=31; D=44566; J=600; K=31000; L=32000; M=31500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1260; Q=25 }
    [PSCustomObject]@{ Row=32; D=44326; J=500; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=33; D=44636; J=500; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=34; D=44656; J=400; K=27000; L=28000; M=27500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1100; Q=25 }
    [PSCustomObject]@{ Row=35; D=44279; J=560; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=36; D=44559; J=360; K=30000; L=31000; M=30500; N="`$/malla 25 kilos"; O="Provincia de Limarí"; P=1220; Q=25 }
    [PSCustomObject]@{ Row=37; D=44559; J=400; K=32000; L=33000; M=32500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1300; Q=25 }
    [PSCustomObject]@{ Row=38; D=44641; J=500; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=39; D=44628; J=560; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=40; D=44582; J=500; K=30000; L=31000; M=30500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1220; Q=25 }
    [PSCustomObject]@{ Row=41; D=44223; J=660; K=32500; L=33000; M=32750; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1310; Q=25 }
    [PSCustomObject]@{ Row=42; D=44505; J=300; K=37000; L=38000; M=37500; N="`$/malla 25 kilos"; O="Perú"; P=1500; Q=25 }
    [PSCustomObject]@{ Row=43; D=44218; J=400; K=34000; L=35000; M=34500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1380; Q=25 }
    [PSCustomObject]@{ Row=44; D=44588; J=500; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=45; D=44239; J=600; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=46; D=44295; J=600; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=47; D=44298; J=500; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=48; D=44643; J=540; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=49; D=44650; J=480; K=27000; L=28000; M=27500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1100; Q=25 }
    [PSCustomObject]@{ Row=50; D=44596; J=500; K=31000; L=32000; M=31500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1260; Q=25 }
    [PSCustomObject]@{ Row=51; D=44307; J=560; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=52; D=44558; J=400; K=30000; L=31000; M=30500; N="`$/malla 25 kilos"; O="Provincia de Limarí"; P=1220; Q=25 }
    [PSCustomObject]@{ Row=53; D=44659; J=400; K=27000; L=28000; M=27500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1100; Q=25 }
    [PSCustomObject]@{ Row=54; D=44333; J=400; K=31500; L=32000; M=31750; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1270; Q=25 }
    [PSCustomObject]@{ Row=55; D=44671; J=540; K=26000; L=27000; M=26500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1060; Q=25 }
    [PSCustomObject]@{ Row=56; D=44309; J=600; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=57; D=44568; J=700; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=58; D=44571; J=600; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=59; D=44238; J=520; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=60; D=44573; J=400; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=61; D=44246; J=600; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=62; D=44251; J=700; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=63; D=44302; J=600; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=64; D=44586; J=600; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=65; D=44260; J=600; K=27000; L=28000; M=27500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1100; Q=25 }
    [PSCustomObject]@{ Row=66; D=44330; J=520; K=32000; L=33000; M=32500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1300; Q=25 }
    [PSCustomObject]@{ Row=67; D=44264; J=400; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=68; D=44253; J=660; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=69; D=44631; J=520; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=70; D=44657; J=500; K=27000; L=28000; M=27500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1100; Q=25 }
    [PSCustomObject]@{ Row=71; D=44335; J=400; K=30000; L=31000; M=30500; N="`$/malla 25 kilos"; O="Provincia de Limarí"; P=1220; Q=25 }
    [PSCustomObject]@{ Row=72; D=44587; J=400; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=73; D=44670; J=600; K=26000; L=27000; M=26500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1060; Q=25 }
    [PSCustomObject]@{ Row=74; D=44664; J=500; K=26000; L=27000; M=26500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1060; Q=25 }
    [PSCustomObject]@{ Row=75; D=44321; J=400; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=76; D=44601; J=300; K=31000; L=32000; M=31500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1260; Q=25 }
    [PSCustomObject]@{ Row=77; D=44265; J=760; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=78; D=44690; J=480; K=24500; L=25000; M=24750; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=990; Q=25 }
    [PSCustomObject]@{ Row=79; D=44209; J=600; K=36000; L=37000; M=36500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1460; Q=25 }
    [PSCustomObject]@{ Row=80; D=44231; J=500; K=25000; L=26000; M=25500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1020; Q=25 }
    [PSCustomObject]@{ Row=81; D=44281; J=640; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=82; D=44316; J=600; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=83; D=44216; J=600; K=36000; L=37000; M=36500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1460; Q=25 }
    [PSCustomObject]@{ Row=84; D=44649; J=400; K=27000; L=28000; M=27500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1100; Q=25 }
    [PSCustomObject]@{ Row=85; D=44678; J=400; K=26000; L=27000; M=26500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1060; Q=25 }
    [PSCustomObject]@{ Row=86; D=44589; J=560; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=87; D=44187; J=400; K=37000; L=38000; M=37500; N="`$/malla 25 kilos"; O="Provincia de Limarí"; P=1500; Q=25 }
    [PSCustomObject]@{ Row=88; D=44629; J=400; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=89; D=44680; J=360; K=27000; L=28000; M=27500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1100; Q=25 }
    [PSCustomObject]@{ Row=90; D=44580; J=600; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=91; D=44594; J=400; K=31000; L=32000; M=31500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1260; Q=25 }
    [PSCustomObject]@{ Row=92; D=44603; J=520; K=31000; L=32000; M=31500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1260; Q=25 }
    [PSCustomObject]@{ Row=93; D=44200; J=400; K=31000; L=32000; M=31500; N="`$/malla 25 kilos"; O="Provincia de Limarí"; P=1260; Q=25 }
    [PSCustomObject]@{ Row=94; D=44278; J=400; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=95; D=44687; J=400; K=25000; L=26000; M=25500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1020; Q=25 }
    [PSCustomObject]@{ Row=96; D=44221; J=460; K=35000; L=36000; M=35500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1420; Q=25 }
    [PSCustomObject]@{ Row=97; D=44293; J=500; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=98; D=44244; J=640; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
    [PSCustomObject]@{ Row=99; D=44245; J=540; K=28000; L=29000; M=28500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1140; Q=25 }
    [PSCustomObject]@{ Row=100; D=44300; J=400; K=29000; L=30000; M=29500; N="`$/malla 25 kilos"; O="Provincia del Elquí"; P=1180; Q=25 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 10).Value = $r.J
    $ws.Cells.Item($r.Row, 11).Value = $r.K
    $ws.Cells.Item($r.Row, 12).Value = $r.L
    $ws.Cells.Item($r.Row, 13).Value = $r.M
    $ws.Cells.Item($r.Row, 14).Value = $r.N
    $ws.Cells.Item($r.Row, 15).Value = $r.O
    $ws.Cells.Item($r.Row, 16).Value = $r.P
    $ws.Cells.Item($r.Row, 17).Value = $r.Q
}

Write-Host "Updated $($rows.Count) rows"